# Insert a new row at position 85 (shifts existing rows 85-195 down to 86-196)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(85).Insert()

$ws.Cells.Item(85, 1).Value2 = 5
$ws.Cells.Item(85, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(85, 3).Value2 = "Maule"
$ws.Cells.Item(85, 4).Value2 = 44915
$ws.Cells.Item(85, 5).Value2 = 7
$ws.Cells.Item(85, 6).Value2 = 100112031
$ws.Cells.Item(85, 7).Value2 = "Poroto verde"
$ws.Cells.Item(85, 8).Value2 = "Sin especificar"
$ws.Cells.Item(85, 9).Value2 = "Primera"
$ws.Cells.Item(85, 10).Value2 = 400
$ws.Cells.Item(85, 11).Value2 = 20000
$ws.Cells.Item(85, 12).Value2 = 20000
$ws.Cells.Item(85, 13).Value2 = 20000
$ws.Cells.Item(85, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(85, 15).Value2 = "Región del Maule"
$ws.Cells.Item(85, 16).Value2 = 800
$ws.Cells.Item(85, 17).Value2 = 25
$ws.Cells.Item(85, 18).Value2 = "Hortaliza"
